# Applies the "Included experimental ROR data and field for delta_g in core models" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new "delta_g" column before the existing "date" column
#    (which was column Q, now becomes R; doi Q->R->S; etc. shift by one)
# ------------------------------------------------------------------
$ws.Columns("Q:Q").Insert()
$ws.Range("Q1").Value = "delta_g"

# ------------------------------------------------------------------
# 2) Backfill the "temperature" (column K) value of 298.15 K for rows
#    that previously had no recorded temperature.
# ------------------------------------------------------------------
$kRows = @(8,9,10,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,47)
foreach ($r in $kRows) {
    $ws.Range("K$r").Value = 298.15
}

# ------------------------------------------------------------------
# 3) Append new experimental ROR (ring-opening redox? / ring-opening
#    NMR) rows 69-86 with the newly measured delta_g values.
# ------------------------------------------------------------------
function Set-RorRow {
    param(
        [int]$Row,
        [string]$MonomerSmiles,
        [double]$DeltaG
    )
    $ws.Range("A$Row").Value = $MonomerSmiles
    $ws.Range("B$Row").Value = "ROR"
    $ws.Range("C$Row").Value = $true
    $ws.Range("E$Row").Value = "CO"
    $ws.Range("H$Row").Value = "CD3OD"
    $ws.Range("I$Row").Value = "s"
    $ws.Range("J$Row").Value = "s"
    $ws.Range("K$Row").Value = 298.15
    $ws.Range("M$Row").Value = "NMR"
    $ws.Range("Q$Row").Value = $DeltaG
    $ws.Range("Q$Row").NumberFormat = "0.00"
    $ws.Range("R$Row").Value = 44950
    $ws.Range("R$Row").NumberFormat = "yyyy\-mm\-dd;@"
    $ws.Range("S$Row").Value = "10.1021/acs.macromol.2c01141"
}

Set-RorRow 69 "O=C1CC(C)O1" -13.137760000000002
Set-RorRow 70 "O=C1OCCCCC1" -10.878400000000001
Set-RorRow 71 "CC1CCCOC1=O" -10.878400000000001
Set-RorRow 72 "O=C([C@H](C)O1)O[C@@H](C)C1=O" -10.878400000000001
Set-RorRow 73 "O=C(CO1)OCC1=O" -10.878400000000001
Set-RorRow 74 "O=C1OCCCC1" -4.0166399999999998
Set-RorRow 75 "O=C1OC(C=C)CCC1CC" -1.7154400000000001
Set-RorRow 76 "O=C1OC(C)CCC1" -0.66944000000000004
Set-RorRow 77 "O=C1OCC(C)CC1" 0.041840000000000002
Set-RorRow 78 "O=C1OC(CCCC)CCC1" 0.041840000000000002
Set-RorRow 79 "O=C1OCCC(C)C1" 1.50624
Set-RorRow 80 "O=C1OC(CC)CCC1CC" 1.50624
Set-RorRow 81 "O=C1OC(C=C)CC/C1=C\C" 2.2175200000000004
Set-RorRow 82 "O=C1OCCC1" 4.8115999999999994
Set-RorRow 83 "O=C1OCC2=C(C=CC=C2)C1" 7.9914399999999999
Set-RorRow 84 "O=C1OC(C)CC1" 9.2466400000000011
Set-RorRow 85 "O=C1OCCC1C" 10.167120000000001
Set-RorRow 86 "O=C1OC(CCCCC)CC1" 10.54368

# ------------------------------------------------------------------
# 4) Cosmetic touches that mirror the author's session: widen the new
#    date column and move the active selection/view.
# ------------------------------------------------------------------
$ws.Columns("R:R").ColumnWidth = 10.15625

$ws.Range("M70").Select()

Write-Output "edit applied"
